$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Modelo" column (F) with a header matching the style of the
# existing header row, then fill in the model name for the single data row.

# Copy the formatting (bold font, border, alignment) of the existing
# "Tipo" header cell (E1) onto the new header cell (F1).
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F1").Value = "Modelo"
$ws.Range("F2").Value = "Pipeline(steps=[('model', LinearRegression())])"
